$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5576.4287
$ws.Range("I113").Value = 4343.222
$ws.Range("J113").Value = 7796.2
$ws.Range("K113").Value = 4343.222
$ws.Range("L113").Value = 7796.2
$ws.Range("M113").Value = -1089.222
$ws.Range("N113").Value = -14304.2
$ws.Range("H116").Value = 2091.077
$ws.Range("I116").Value = 1997.7142
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 1997.7142
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1444.2858
$ws.Range("N116").Value = -9084
$ws.Range("H137").Value = 1438.081
$ws.Range("I137").Value = 920.15
$ws.Range("J137").Value = 2047.4117
$ws.Range("K137").Value = 2760.45
$ws.Range("L137").Value = 6142.2351
$ws.Range("M137").Value = -210.4499999999998
$ws.Range("N137").Value = -11242.2351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1609.6
$ws.Range("I74").Value = 1302.9231
$ws.Range("J74").Value = 2179.1428
$ws.Range("K74").Value = 1302.9231
$ws.Range("L74").Value = 2179.1428
$ws.Range("M74").Value = -428.9231
$ws.Range("N74").Value = -3927.1428
$ws.Range("H77").Value = 1609.6
$ws.Range("I77").Value = 1302.9231
$ws.Range("J77").Value = 2179.1428
$ws.Range("K77").Value = 6514.6155
$ws.Range("L77").Value = 10895.714
$ws.Range("M77").Value = -2146.6155
$ws.Range("N77").Value = -19631.714
$ws.Range("H122").Value = 25641024
$ws.Range("I122").Value = 25641024
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 76923072
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -76920622
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1375
$ws.Range("I64").Value = 1250
$ws.Range("J64").Value = 1500
$ws.Range("K64").Value = 1250
$ws.Range("L64").Value = 1500
$ws.Range("M64").Value = -1025
$ws.Range("N64").Value = -1950
$ws.Range("H67").Value = 1375
$ws.Range("I67").Value = 1250
$ws.Range("J67").Value = 1500
$ws.Range("K67").Value = 1250
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = -470
$ws.Range("N67").Value = -3060
$ws.Range("H80").Value = 545.8570999999999
$ws.Range("J80").Value = 719.5
$ws.Range("L80").Value = 719.5
$ws.Range("N80").Value = -2715.5
$ws.Range("H83").Value = 545.8570999999999
$ws.Range("J83").Value = 719.5
$ws.Range("L83").Value = 3597.5
$ws.Range("N83").Value = -13581.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1529.25
$ws.Range("I58").Value = 1349.0625
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 1349.0625
$ws.Range("L58").Value = 2250
$ws.Range("M58").Value = -1146.0625
$ws.Range("N58").Value = -2656
$ws.Range("H68").Value = 30500
$ws.Range("J68").Value = 30500
$ws.Range("L68").Value = 30500
$ws.Range("N68").Value = -31998
$ws.Range("H71").Value = 30500
$ws.Range("J71").Value = 30500
$ws.Range("L71").Value = 91500
$ws.Range("N71").Value = -98988
$ws.Range("H99").Value = 12523790
$ws.Range("I99").Value = 30000
$ws.Range("J99").Value = 17878272
$ws.Range("K99").Value = 30000
$ws.Range("L99").Value = 17878272
$ws.Range("M99").Value = -28502
$ws.Range("N99").Value = -17881268
$ws.Range("H126").Value = 12523790
$ws.Range("I126").Value = 30000
$ws.Range("J126").Value = 17878272
$ws.Range("K126").Value = 90000
$ws.Range("L126").Value = 53634816
$ws.Range("M126").Value = -87530
$ws.Range("N126").Value = -53639756
$ws.Range("H136").Value = 1529.25
$ws.Range("I136").Value = 1349.0625
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 4047.1875
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -1497.1875
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9669.833000000001
$ws.Range("J5").Value = 2400.8333
$ws.Range("L5").Value = 7202.499899999999
$ws.Range("N5").Value = -7426.499899999999
$ws.Range("H34").Value = 1709.8334
$ws.Range("I34").Value = 547.7692
$ws.Range("J34").Value = 2366.652
$ws.Range("K34").Value = 1643.3076
$ws.Range("L34").Value = 7099.956
$ws.Range("M34").Value = -1559.3076
$ws.Range("N34").Value = -7267.956
$ws.Range("H135").Value = 9669.833000000001
$ws.Range("J135").Value = 2400.8333
$ws.Range("L135").Value = 21607.4997
$ws.Range("N135").Value = -26677.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23922.111
$ws.Range("J15").Value = 23922.111
$ws.Range("L15").Value = 23922.111
$ws.Range("N15").Value = -24498.111
$ws.Range("H81").Value = 23922.111
$ws.Range("J81").Value = 23922.111
$ws.Range("L81").Value = 23922.111
$ws.Range("N81").Value = -25918.111
$ws.Range("H84").Value = 23922.111
$ws.Range("J84").Value = 23922.111
$ws.Range("L84").Value = 71766.333
$ws.Range("N84").Value = -81750.333
$ws.Range("H122").Value = 3243151.8
$ws.Range("I122").Value = 4323229
$ws.Range("J122").Value = 2919.6
$ws.Range("K122").Value = 12969687
$ws.Range("L122").Value = 8758.799999999999
$ws.Range("M122").Value = -12967237
$ws.Range("N122").Value = -13658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5091616
$ws.Range("I122").Value = 5104490
$ws.Range("J122").Value = 5001500
$ws.Range("K122").Value = 15313470
$ws.Range("L122").Value = 15004500
$ws.Range("M122").Value = -15311020
$ws.Range("N122").Value = -15009400
$ws.Range("H132").Value = 17338158
$ws.Range("I132").Value = 19701794
$ws.Range("J132").Value = 4833.3335
$ws.Range("K132").Value = 59105382
$ws.Range("L132").Value = 14500.0005
$ws.Range("M132").Value = -59102852
$ws.Range("N132").Value = -19560.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2051.3914
$ws.Range("I132").Value = 2051.3914
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6154.174199999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3624.174199999999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1234.625
$ws.Range("I136").Value = 1042.6154
$ws.Range("J136").Value = 2066.6667
$ws.Range("K136").Value = 3127.8462
$ws.Range("L136").Value = 6200.000100000001
$ws.Range("M136").Value = -577.8462
$ws.Range("N136").Value = -11300.0001
